# Insert a new data row above row 139 (shifts existing rows 139:171 down to 140:172)
# and populate it with the new daily price-report record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(139).Insert()

$ws.Cells.Item(139, 1).Value = 10
$ws.Cells.Item(139, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(139, 3).Value = "La Araucanía"
$ws.Cells.Item(139, 4).Value = 44508
$ws.Cells.Item(139, 5).Value = 9
$ws.Cells.Item(139, 6).Value = 100112039
$ws.Cells.Item(139, 7).Value = "Ciboulette"
$ws.Cells.Item(139, 8).Value = "Sin especificar"
$ws.Cells.Item(139, 9).Value = "Primera"
$ws.Cells.Item(139, 10).Value = 40
$ws.Cells.Item(139, 11).Value = 5000
$ws.Cells.Item(139, 12).Value = 5000
$ws.Cells.Item(139, 13).Value = 5000
$ws.Cells.Item(139, 14).Value = "$/docena de atados"
$ws.Cells.Item(139, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(139, 16).Value = 1667
$ws.Cells.Item(139, 17).Value = 3
$ws.Cells.Item(139, 18).Value = "Hortaliza"
